$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the "Organization" block (rows 30:32, columns B:H) down to rows 33:35 to
# create the new "Person" block with identical formatting.
$src = $ws.Range("B30:H32")
$dst = $ws.Range("B33:H35")
$src.Copy($dst)

# The copy carries style 19 onto C35; the target uses style 20 (no top border), same as C33/C34.
$ws.Range("C35").Borders.Item(8).LineStyle = -4142

# Row 33: Person - Default access
$ws.Range("C33").Value2 = "PERSON"
$ws.Range("B33").Value2 = "Person – Default access"
$ws.Range("D33").Value2 = "participants.?[participantType == '*'].isEmpty()"
$ws.Range("G33").Value2 = "*, *"

# Row 34: Person - Default owner
$ws.Range("B34").Value2 = "Person – Default owner"
$ws.Range("C34").Value2 = "PERSON"
$ws.Range("D34").Value2 = "participants.?[participantType == 'owner'].isEmpty()"
$ws.Range("H34").Value2 = "owner, creator"

# Row 35: Person - Default group
$ws.Range("B35").Value2 = "Person – Default group"
$ws.Range("C35").Value2 = "PERSON"
$ws.Range("D35").Value2 = "participants.?[participantType == 'owning group'].isEmpty()"
$ws.Range("G35").Value2 = "owning group, ACM_ADMINISTRATOR_DEV"

# Update the view: scroll so row 27 is at the top and select B36 (the cell just below the new data)
$ws.Range("B36").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 27
$win.ScrollColumn = 1

Write-Host "done"
